$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$2b$12$dnEm.7Uhx87Hlj5HIMIV4OtKrgEGbpyMzKxk1Fj.05svPQcCe3vB.'
$ws.Range("B3").Value = '$2b$12$iNQxdTs2Q2pdmYyYkvYrke2tRsiVaHHag8T/LIBqdp/kT5PldprVK'

$ws.Range("A4").Value = "professor"
$ws.Range("B4").Value = '$2b$12$URKLT7SIDakWrHc7tUd2jOeDZdk/aIkTFKcBrJ6PG7TacpI/hQnFy'

$ws.Range("A5").Value = "aluno"
$ws.Range("B5").Value = '$2b$12$PERvjXXRGAkqcChk99atF.pmHy/l//vu9z5En1N5ovoW6qsE28UOG'
